$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text number format to all cells we are about to write,
# so Excel does not auto-convert numeric-looking strings (e.g. "22.09",
# "0.000009464") into actual numbers/scientific notation.
$targetCells = @(
    'D2', 'E2', 'D3', 'E3', 'E4', 'D5', 'E5', 'E6', 'D7', 'E7', 'D8', 'E8', 'D9', 'E9', 'D10', 'E10', 'D11', 'E11', 'D12', 'E12', 'D13', 'E13', 'D14', 'E14', 'D15', 'E15', 'E16', 'D17', 'E17', 'D18', 'E18', 'D19', 'E19', 'E20', 'D21', 'E21', 'D22', 'E22', 'D23', 'E23', 'D24', 'E24', 'D25', 'E25', 'D26', 'E26', 'D27', 'E27', 'D28', 'E28', 'E29', 'D30', 'E30', 'D31', 'E31', 'D32', 'E32', 'D33', 'E33', 'D34', 'E34', 'D35', 'E35', 'D36', 'E36', 'E37', 'D38', 'E38', 'D39', 'E39', 'D40', 'E40', 'D41', 'E41', 'D42', 'E42', 'D43', 'E43', 'D44', 'E44', 'D45', 'E45', 'E46', 'D47', 'E47', 'B48', 'C48', 'D48', 'E48', 'B49', 'C49', 'D49', 'E49', 'B50', 'C50', 'D50', 'E50', 'D51', 'E51'
)
foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the new values cell by cell
$ws.Range('D2').Value = '28.912.18'
$ws.Range('E2').Value = '  -1.63%  '
$ws.Range('D3').Value = '1.910.46'
$ws.Range('E3').Value = '  -1.71%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '325.01'
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('D7').Value = '0.4589'
$ws.Range('E7').Value = '  -0.87%  '
$ws.Range('D8').Value = '0.3815'
$ws.Range('E8').Value = '  -1.51%  '
$ws.Range('D9').Value = '0.07719'
$ws.Range('E9').Value = '  -1.46%  '
$ws.Range('D10').Value = '0.9798'
$ws.Range('E10').Value = '  +0.43%  '
$ws.Range('D11').Value = '22.09'
$ws.Range('E11').Value = '  -2.45%  '
$ws.Range('D12').Value = '1.918.66'
$ws.Range('E12').Value = '  -1.34%  '
$ws.Range('D13').Value = '6.933'
$ws.Range('E13').Value = '  -2.24%  '
$ws.Range('D14').Value = '5.660'
$ws.Range('E14').Value = '  -1.77%  '
$ws.Range('D15').Value = '0.07043'
$ws.Range('E15').Value = '  -0.18%  '
$ws.Range('E16').Value = '  +0.00%  '
$ws.Range('D17').Value = '83.80'
$ws.Range('E17').Value = '  -3.43%  '
$ws.Range('D18').Value = '0.000009464'
$ws.Range('E18').Value = '  -3.70%  '
$ws.Range('D19').Value = '16.66'
$ws.Range('E19').Value = '  -2.54%  '
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('D21').Value = '28.903.91'
$ws.Range('E21').Value = '  -1.73%  '
$ws.Range('D22').Value = '5.321'
$ws.Range('E22').Value = '  -2.86%  '
$ws.Range('D23').Value = '10.88'
$ws.Range('E23').Value = '  -1.60%  '
$ws.Range('D24').Value = '2.095'
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').Value = '158.36'
$ws.Range('E25').Value = '  +0.89%  '
$ws.Range('D26').Value = '19.05'
$ws.Range('E26').Value = '  -1.71%  '
$ws.Range('D27').Value = '5.655'
$ws.Range('E27').Value = '  -1.71%  '
$ws.Range('D28').Value = '117.74'
$ws.Range('E28').Value = '  -0.62%  '
$ws.Range('E29').Value = '  +0.36%  '
$ws.Range('D30').Value = '0.09294'
$ws.Range('E30').Value = '  -0.55%  '
$ws.Range('D31').Value = '0.8668'
$ws.Range('E31').Value = '  +0.60%  '
$ws.Range('D32').Value = '5.077'
$ws.Range('E32').Value = '  -2.00%  '
$ws.Range('D33').Value = '1.250'
$ws.Range('E33').Value = '  -4.35%  '
$ws.Range('D34').Value = '3.124'
$ws.Range('E34').Value = '  +0.62%  '
$ws.Range('D35').Value = '0.05715'
$ws.Range('E35').Value = '  -1.07%  '
$ws.Range('D36').Value = '1.161'
$ws.Range('E36').Value = '  +0.59%  '
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('D38').Value = '0.02046'
$ws.Range('E38').Value = '  -1.91%  '
$ws.Range('D39').Value = '7.409'
$ws.Range('E39').Value = '  -3.78%  '
$ws.Range('D40').Value = '0.5488'
$ws.Range('E40').Value = '  -3.35%  '
$ws.Range('D41').Value = '0.1755'
$ws.Range('E41').Value = '  -1.52%  '
$ws.Range('D42').Value = '2.867'
$ws.Range('E42').Value = '  +4.82%  '
$ws.Range('D43').Value = '9.309'
$ws.Range('E43').Value = '  -1.13%  '
$ws.Range('D44').Value = '0.5174'
$ws.Range('E44').Value = '  -2.30%  '
$ws.Range('D45').Value = '11.21'
$ws.Range('E45').Value = '  -1.97%  '
$ws.Range('E46').Value = '  +0.53%  '
$ws.Range('D47').Value = '2.088'
$ws.Range('E47').Value = '  +0.31%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '1.779'
$ws.Range('E48').Value = '  -2.11%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = '110.30'
$ws.Range('E49').Value = '  -0.94%  '
$ws.Range('B50').Value = 'PEPE'
$ws.Range('C50').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D50').Value = '0.000002557'
$ws.Range('E50').Value = '  -9.16%  '
$ws.Range('D51').Value = '0.2875'
$ws.Range('E51').Value = '  -4.20%  '
